$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated D/E measurements per block (row -> D,E) - "no exch" run results
$updates = @(
    @{Row=4; D=1.6885399999999999; E=0.08319},
    @{Row=5; D=1.6885399999999999; E=0.072370000000000004},
    @{Row=6; D=1.6885399999999999; E=0.069309999999999997},
    @{Row=7; D=1.6885399999999999; E=0.070620000000000002},
    @{Row=8; D=1.6885399999999999; E=0.068409999999999999},
    @{Row=9; D=1.6885399999999999; E=0.070680000000000007},
    @{Row=11; D=3.9244400000000002; E=0.16545000000000001},
    @{Row=12; D=3.9244400000000002; E=0.14341000000000001},
    @{Row=13; D=3.9244400000000002; E=0.13658999999999999},
    @{Row=14; D=3.9244400000000002; E=0.14051},
    @{Row=15; D=3.9244400000000002; E=0.13414999999999999},
    @{Row=16; D=3.9244400000000002; E=0.14394000000000001},
    @{Row=18; D=8.83385; E=0.35424},
    @{Row=19; D=8.83385; E=0.30032999999999999},
    @{Row=20; D=8.83385; E=0.29984},
    @{Row=21; D=8.83385; E=0.30203000000000002},
    @{Row=22; D=8.83385; E=0.30414000000000002},
    @{Row=23; D=8.83385; E=0.30220000000000002},
    @{Row=25; D=15.68573; E=0.67940999999999996},
    @{Row=26; D=15.68573; E=0.62417},
    @{Row=27; D=15.68573; E=0.64097000000000004},
    @{Row=28; D=15.68573; E=0.64605999999999997},
    @{Row=29; D=15.68573; E=0.64949000000000001},
    @{Row=30; D=15.68573; E=0.65795999999999999},
    @{Row=35; D=0.89002000000000003; E=0.080299999999999996},
    @{Row=36; D=0.89002000000000003; E=0.072419999999999998},
    @{Row=37; D=0.89002000000000003; E=0.066689999999999999},
    @{Row=38; D=0.89002000000000003; E=0.066629999999999995},
    @{Row=39; D=0.89002000000000003; E=0.065119999999999997},
    @{Row=40; D=0.89002000000000003; E=0.077270000000000005},
    @{Row=42; D=1.9026099999999999; E=0.16302},
    @{Row=43; D=1.9026099999999999; E=0.12903000000000001},
    @{Row=44; D=1.9026099999999999; E=0.12418999999999999},
    @{Row=45; D=1.9026099999999999; E=0.12834000000000001},
    @{Row=46; D=1.9026099999999999; E=0.12559000000000001},
    @{Row=47; D=1.9026099999999999; E=0.14524000000000001},
    @{Row=49; D=4.5672499999999996; E=0.29285},
    @{Row=50; D=4.5672499999999996; E=0.24917},
    @{Row=51; D=4.5672499999999996; E=0.22931000000000001},
    @{Row=52; D=4.5672499999999996; E=0.22932},
    @{Row=53; D=4.5672499999999996; E=0.23227},
    @{Row=54; D=4.5672499999999996; E=0.28586},
    @{Row=56; D=10.55977; E=0.57065999999999995},
    @{Row=57; D=10.55977; E=0.48305999999999999},
    @{Row=58; D=10.55977; E=0.45971000000000001},
    @{Row=59; D=10.55977; E=0.46072999999999997},
    @{Row=60; D=10.55977; E=0.45900999999999997},
    @{Row=61; D=10.55977; E=0.55408000000000002}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Update the selection to match the recorded view state after the edit
[void]$ws.Range("D56:D61").Select()
